$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for the new "Save" column, matching the bold/bordered style
# already used by the other header cells (e.g. G1).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Save values per row (1 = save recorded, 0 = no save), row 2 through 14.
$saveValues = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 0
    6  = 1
    7  = 1
    8  = 0
    9  = 1
    10 = 0
    11 = 0
    12 = 0
    13 = 1
    14 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
